# Apply updates to column F (dSF) values as described by the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value  = -6
$ws.Range("F5").Value  = -1
$ws.Range("F6").Value  = -8
$ws.Range("F10").Value = 3
$ws.Range("F14").Value = -1
$ws.Range("F17").Value = -1
$ws.Range("F19").Value = -6
$ws.Range("F28").Value = 1
$ws.Range("F29").Value = -3
